{"js": "// The document contains a single 5-column table whose column widths need a\n// small adjustment (in twentieths-of-a-point == dxa/twips):\n//   column 1 (index 0): 714 dxa -> 712 dxa  (35.7pt -> 35.6pt)\n//   column 3 (index 2): 1991 dxa -> 1993 dxa (99.55pt -> 99.65pt)\n// Word JS `TableCell.columnWidth` is column-wide: setting it on any cell in\n// a column resizes every cell in that column (and the table's <w:gridCol>).\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// dxa -> points: Office.js widths are expressed in points, OOXML widths\n// (w:w under w:tcW / w:gridCol) are in twentieths of a point (dxa).\nconst dxaToPoints = (dxa) => dxa / 20;\n\n// First column: 714 dxa -> 712 dxa.\ntable.getCell(0, 0).columnWidth = dxaToPoints(712);\n// Third column: 1991 dxa -> 1993 dxa.\ntable.getCell(0, 2).columnWidth = dxaToPoints(1993);\n\nawait context.sync();\n", "ps1": "# The document contains a single 5-column table whose column widths need a\n# small adjustment (values below are in twentieths-of-a-point == dxa/twips,\n# matching the <w:tcW>/<w:gridCol> w:w attributes in the OOXML):\n#   column 1: 714 dxa -> 712 dxa  (35.7pt  -> 35.6pt)\n#   column 3: 1991 dxa -> 1993 dxa (99.55pt -> 99.65pt)\n# Word's Table.Columns(n).Width is column-wide: setting it resizes every\n# cell in that column (and the table's <w:gridCol> grid definition) at once.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# dxa -> points: COM widths are expressed in points, OOXML widths\n# (w:w under w:tcW / w:gridCol) are in twentieths of a point (dxa).\n$table.Columns.Item(1).Width = 712 / 20\n$table.Columns.Item(3).Width = 1993 / 20\n"}
